# TC_5-FDR_E2E-FDR-2938-SYSTEM ADMIN REFERENCE DATA
#
# Update the draw-date seed values on the single data sheet
# (FDR_End_End_Reject_TestData): push drawDate (col A, rows 2-3) forward
# from 9/30/2020 to 10/12/2020. payDate (col B) is a dependent formula
# (=A+2) and recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# drawDate: 9/30/2020 (44104) -> 10/12/2020 (44116) for both data rows
$ws.Range("A2").Value = 44116
$ws.Range("A3").Value = 44116

# Column A was widened slightly to fit the new content (previously shared
# the same width as column B); give it its own, slightly wider column
# definition.
$ws.Columns.Item(1).ColumnWidth = 11

# Move the active selection from B4 to A4
$ws.Range("A4").Select()
